# Add data for 2022-11-30: update the sheet name, header label, and
# November/Total values for the "Total" column (I).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab to reflect the new "through" date.
$ws.Name = "Through 2022-11-22"

# Update the header label in I1 (shared string "2022 (through 11-21)" -> "... 11-22").
$ws.Range("I1").Value = "2022 (through 11-22)"

# Update November (row 12) and Total (row 14) figures in the Total column (I).
$ws.Range("I12").Value = 83
$ws.Range("I14").Value = 1480
